# Jade Skins powers definition
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("powerups")

# ---------------------------------------------------------------------------
# 1. Insert 4 new rows just below the existing table (old rows 165-168),
#    pushing the trailing "FREEZE CONSTANT" block from 165-169 -> 169-173.
# ---------------------------------------------------------------------------
$ws.Rows("165:168").Insert()

# ---------------------------------------------------------------------------
# 2. Fix up pre-existing formatting glitches shown by the diff.
# ---------------------------------------------------------------------------
# B90 loses its bold font (style 80 -> style 21, i.e. same look as B91/B92...)
$ws.Range("B91").Copy()
$ws.Range("B90").PasteSpecial(-4122) | Out-Null
$ws.Range("B90").Value = "disguise_transform_gold_LOWx2_damage_arrows"

# Row 120 drops its thick-bottom divider formatting.
$ws.Range("A119").Copy()
$ws.Range("A120").PasteSpecial(-4122) | Out-Null
$ws.Rows("120:120").RowHeight = 15.75

# ---------------------------------------------------------------------------
# 3. Populate the two plain new data rows (165 & 166) - same look as row 164.
# ---------------------------------------------------------------------------
$ws.Range("A164:L164").Copy()
$ws.Range("A165:L165").PasteSpecial(-4122) | Out-Null
$ws.Range("A164:L164").Copy()
$ws.Range("A166:L166").PasteSpecial(-4122) | Out-Null

$ws.Range("A165").Value = "<Definition>"
$ws.Range("B165").Value = "disguise_boost_LOW_faster_boost"
$ws.Range("C165").Value = "combined"
$ws.Range("D165").Value = "stats"
$ws.Range("E165").Value = "disguise_boost_LOW"
$ws.Range("F165").Value = "disguise_faster_boost"
$ws.Range("G165").Value = "icon_power_special"
$ws.Range("H165").Value = "icon_special"
$ws.Range("I165").Value = "TID_POWERUP_DISGUISE_51_NAME"
$ws.Range("J165").Value = "TID_POWERUP_DISGUISE_51_DESC"
$ws.Range("K165").Value = "TID_POWERUP_DISGUISE_51_DESC_SHORT"
$ws.Range("L165").Value = 2

$ws.Range("A166").Value = "<Definition>"
$ws.Range("B166").Value = "disguise_speed_fury_duration_LOW"
$ws.Range("C166").Value = "combined"
$ws.Range("D166").Value = "stats"
$ws.Range("E166").Value = "disguise_speed"
$ws.Range("F166").Value = "disguise_fury_duration_LOW"
$ws.Range("G166").Value = "icon_power_special"
$ws.Range("H166").Value = "icon_special"
$ws.Range("I166").Value = "TID_POWERUP_DISGUISE_52_NAME"
$ws.Range("J166").Value = "TID_POWERUP_DISGUISE_52_DESC"
$ws.Range("K166").Value = "TID_POWERUP_DISGUISE_52_DESC_SHORT"
$ws.Range("L166").Value = 2

# ---------------------------------------------------------------------------
# 4. Row 167 is the closing row of the group - "end of block" styling.
# ---------------------------------------------------------------------------
$ws.Range("A14:G14").Copy()
$ws.Range("A167:G167").PasteSpecial(-4122) | Out-Null
$ws.Range("H23").Copy()
$ws.Range("H167").PasteSpecial(-4122) | Out-Null
$ws.Range("I14:K14").Copy()
$ws.Range("I167:K167").PasteSpecial(-4122) | Out-Null
$ws.Range("G67").Copy()
$ws.Range("L167").PasteSpecial(-4122) | Out-Null

$ws.Range("A167").Value = "<Definition>"
$ws.Range("B167").Value = "disguise_speed_LOW_boost_LOW_free_revive"
$ws.Range("C167").Value = "combined"
$ws.Range("D167").Value = "stats"
$ws.Range("E167").Value = "disguise_speed_LOW_boost_LOW"
$ws.Range("F167").Value = "free_revive"
$ws.Range("G167").Value = "icon_power_special"
$ws.Range("H167").Value = "icon_special"
$ws.Range("I167").Value = "TID_POWERUP_DISGUISE_53_NAME"
$ws.Range("J167").Value = "TID_POWERUP_DISGUISE_53_DESC"
$ws.Range("K167").Value = "TID_POWERUP_DISGUISE_53_DESC_SHORT"
$ws.Range("L167").Value = 3

# Row-level "bottom divider" formatting for the new last row of the group.
$ws.Rows("167:167").Borders.Item(9).LineStyle = 1
$ws.Rows("167:167").Borders.Item(9).Weight = 2

# ---------------------------------------------------------------------------
# 5. Grow the table to cover the 3 new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:L167"))

Write-Host "done"
